$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.875.62"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.303.68"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.71"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -4.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.292.85"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.185"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -8.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.583"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.63"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.17%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.65"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.06%  "

$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "634.44"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.827.48"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.15"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.677.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.83%  "

$ws.Range("E19").Value = "  -2.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.301.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.40"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.906"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.11"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.38%  "

$ws.Range("E26").Value = "  -6.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.60"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.68"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.22"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.38"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("E34").Value = "  -2.83%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "544.20"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.17%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.105"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.87%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.818.09"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.48"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.42%  "

$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0740"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "34.02"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.54%  "

$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.49%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.25"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.63%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.129"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.73%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.336"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -8.34%  "

$ws.Range("B46").Value = "CoreDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.21"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -15.73%  "

$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.66%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.128"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.66%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.15%  "
